# Added inductors to component list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 - Inductor Bat cntr (Bourns)
$ws.Range("A35").Value = "Inductor Bat cntr"
$ws.Range("B35").Value = "https://www.mouser.ee/ProductDetail/Bourns/SRP1265C-4R7M?qs=OlC7AqGiEDkg4Xkb2TDlFw%3D%3D"
$ws.Range("C35").Value = 1.5
$ws.Range("D35").Value = 1

# Row 36 - Inductor 5V (Walsin)
$ws.Range("A36").Value = "Inductor 5V"
$ws.Range("B36").Value = "https://www.mouser.ee/ProductDetail/Walsin/WLPMA0A040M4R7LC?qs=B6kkDfuK7%2FAGhXqtZ1HOFg%3D%3D"
$ws.Range("C36").Value = 0.67
$ws.Range("D36").Value = 1

# Row 37 - Inductor USB source (Taiyo-Yuden)
# Note: link (B37) is entered before the label (A37) so new shared-string
# indices land in the same order as the source workbook.
$ws.Range("B37").Value = "https://www.mouser.ee/ProductDetail/Taiyo-Yuden/NRS8030T1R0NJGJ?qs=PzICbMaShUfBDq1Kfb1D%252Bg%3D%3D"
$ws.Range("A37").Value = "Inductor USB source"
$ws.Range("C37").Value = 0.41
$ws.Range("D37").Value = 2

# Update the view/selection to match where the new rows were edited.
$ws.Range("F36").Select() | Out-Null
